$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the inverted upper/lower boundary for the methanotrophic source (row 3)
$ws.Range("B3").Value = -120
$ws.Range("C3").Value = -40

# Update the active selection to B4
$ws.Range("B4").Select()
